# Re-generate the list of files in column A (the review dates in column B stay fixed per row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$files = @(
    'Routine_Care/Nursing for Arterial and Central Venous Lines.pdf',
    'Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf',
    'Routine_Care/Faecal  incontinence skin care.pdf',
    'Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf',
    'Drugs/heparin_critical_care_only.pdf',
    'Routine_Care/Invasive Flush Systems.pdf',
    'GI_Liver_and_Transplant/Pancreatic Irrigation.pdf',
    'Neurological/Thiopentone levels.pdf',
    'Infection_and_sepsis/Ebola/Ebola.pdf',
    'GI_Liver_and_Transplant/Nasal bridle.pdf',
    'Breathing(Respiratory)/Equipment/AMBU AScope.pdf',
    'Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf',
    'Routine_Care/Central venous catheter removal.pdf',
    'Routine_Care/Tracheostomy_nursing_care.pdf',
    'Airway/Critical care extubation checklist.pdf',
    'ECLS/Extra Corporeal Carbon Dioxide Removal.pdf',
    'Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf',
    'Trauma and Burns/Mangement of burns.pdf',
    'Drugs/diazepam_diazemuls.pdf',
    'Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf',
    'End_of_life_care/Reasons to report a death to PF.pdf',
    'Neurological/SOP -  Femoral site care.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf',
    'Drugs/sodium_bicarbonate.pdf',
    'Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf',
    'Post_op_care/Epidural top-up.pdf',
    'Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf',
    'Drugs/heparin for Haemofiltration.pdf',
    'Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf',
    'Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf',
    'Covid-19/WGH/CoVid intubation checklist WGH.pdf',
    'Airway/Emergency intubation checklist_em_pub.pdf',
    'Drugs/fentanyl.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf',
    'Covid-19/WGH/WGH_CT_Transfer_May.pdf',
    'Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf',
    'Organ_donation/Donation after circulatory death.pdf',
    'Airway/Percutaneous tracheostomy checklist.pdf',
    'Delirium/Managing a Potentially Violent Patient.pdf',
    'Delirium/Risk assessment posi mit.pdf',
    'Infection_and_sepsis/SOP Ultrasound Cleaning.pdf',
    'GI_Liver_and_Transplant/Treatment of constipation.pdf',
    'GI_Liver_and_Transplant/Abdominal pressure measurement.pdf',
    'Airway/Anticipated difficult airway tool.pdf',
    'Drugs/ketamine_in_asthma.pdf',
    'Breathing(Respiratory)/HFNO.pdf',
    'Airway/McGrath Mac.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf',
    'Delirium/Drugs Causing Delirium and Agitiation.pdf',
    'Neurological/Sub arachnoid haemorrhage management.pdf',
    'End_of_life_care/Documentation following death.pdf',
    'Drugs/zanamivir.pdf',
    'Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf',
    'Drugs/insulin.pdf',
    'Breathing(Respiratory)/Equipment/HFNO Set Up.pdf',
    'Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf',
    'Cardiovascular/Steroids for Septic Shock.pdf',
    'Breathing(Respiratory)/Equipment/APRV.pdf',
    'Post_op_care/Epidural Haematoma.pdf',
    'Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf',
    'Breathing(Respiratory)/Equipment/T piece Y piece.pdf',
    'Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf',
    'Drugs/midazolam and thiopental levels.pdf',
    'Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf',
    'Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf',
    'Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf',
    'Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf',
    'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf',
    'Drugs/Antibiotic doses in CVVHD.pdf',
    'Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf',
    'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf',
    'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf',
    'ECLS/RIE ECLS Anti Xa Protocol.pdf',
    'Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf',
    'Drugs/vasopressin_sepsis.pdf',
    'Drugs/vasopressin organ donation.pdf',
    'Covid-19/videos/Donning and Doffing Video.pdf',
    'Transfer/ACCP Transfers.pdf',
    'Breathing(Respiratory)/Equipment/Bipap V60.pdf',
    'Breathing(Respiratory)/CPAP.pdf',
    'Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf',
    'Infection_and_sepsis/Infection indications for IVIG.pdf',
    'Drugs/piperacillin_tazobactam extended_infusion.pdf',
    'Procedures/CVC Guidance/Securing CVCs.pdf',
    'Covid-19/Covid 19 Death Certification Guideline.pdf',
    'Neurological/Treatment of status epilepticus.pdf',
    'Routine_Care/Video Communication.pdf',
    'Cardiovascular/Cardiogenic Shock.pdf',
    'Drugs/isoprenaline.pdf',
    'Haematology_CAR-T/Haem_ICU_transfer.pdf',
    'Cardiovascular/Management of hypertension within Critical Care.pdf',
    'Drugs/aminophylline.pdf',
    'Haematology_CAR-T/CRS.pdf',
    'Drugs/pancuronium.pdf',
    'Haematology_CAR-T/ICANS.pdf',
    'Drugs/phenytoin.pdf',
    'Drugs/rocuronium.pdf',
    'Drugs/milrinone.pdf',
    'Policies_and_admin/General Critical Care SOP_pub.pdf',
    'GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf',
    'Drugs/clonidine.pdf',
    'Neurological/Critical Care MRI Procedure_pub.pdf',
    'Ethics_and_Law/DNACPR policy for Scotland.pdf',
    'End_of_life_care/Palliative extubation & withdrawal of invasive ventilatory support nursing checklist.pdf',
    'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf',
    'Neurological/Management of traumatic brain injury.pdf',
    'Organ_donation/Organ Retrieval SOP.pdf',
    'Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf',
    'Infection_and_sepsis/Trip Out of Unit infection guidance.pdf',
    'Ethics_and_Law/Care at the End of Life (FICM).pdf',
    'Airway/Cook Staged Extubation Set.pdf',
    'Drugs/noradrenaline (central).pdf',
    'Post_op_care/Epidural hypotension.pdf',
    'Drugs/glyceryl_trinitrate.pdf',
    'GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf',
    'GI_Liver_and_Transplant/Fulminant Liver Failure.pdf',
    'Drugs/dexmedetomidine.pdf',
    'Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf',
    'Infection_and_sepsis/Antifungal guidance in critical care.pdf',
    'Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf',
    'Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf',
    'Transfer/Transfer Outdoors to Garden Guideline.pdf',
    'Breathing(Respiratory)/ARDS Strategy.pdf',
    'Drugs/dobutamine.pdf',
    'Drugs/adrenaline.pdf',
    'Drugs/hydralazine.pdf',
    'Post_op_care/Major OMFS Free Flap.pdf',
    'Drugs/Alteplase for massive PE.pdf',
    'Drugs/alfentanil.pdf',
    'Drugs/magnesium.pdf',
    'Drugs/all IV drug infusion information.pdf',
    'Drugs/neostigmine.pdf',
    'Drugs/vancomycin.pdf',
    'Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf',
    'Drugs/labetalol.pdf',
    'Neurological/Intrathecal policy RIE.pdf',
    'Drugs/midazolam.pdf',
    'Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf',
    'Drugs/potassium.pdf',
    'Procedures/CVC Guidance/CVC NHL  April 2023.pdf',
    'Drugs/salbutamol.pdf',
    'Drugs/nimodipine.pdf',
    'Routine_Care/ICU Eye Care Guideline.pdf',
    'Drugs/amiodarone.pdf',
    'Drugs/nicardipine.pdf',
    'Drugs/phenobarbitone.pdf',
    'Procedures/Arterial Line insertion for ACCPs.pdf',
    'Breathing(Respiratory)/Manual Ventilation and MHI.pdf',
    'Drugs/noradrenaline (peripheral).pdf',
    'Neurological/Ventriculitis Guideline.pdf',
    'Drugs/morphine.pdf',
    'Drugs/Epoprostenol.pdf',
    'Drugs/phenylephrine.pdf',
    'Cardiovascular/Cardiac Output Monitoring_pub .pdf',
    'Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf',
    'Drugs/calcium.pdf',
    'Drugs/dalteparin_thromboprophylaxis.pdf',
    'Drugs/Vancomycin Continuous Infusion Fluid Restricted.pdf',
    'Drugs/atracurium.pdf',
    'Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf',
    'Drugs/valproate.pdf',
    'Policies_and_admin/Pet Visitation.pdf',
    'Procedures/Inadvertent Catheter Placement Guideline.pdf',
    'Drugs/stress ulcer prophylaxis.pdf',
    'Drugs/ketamine_for_status epilepticus.pdf',
    'Drugs/Phosphate.pdf',
    'Drugs/Thiopentone.pdf',
    'Breathing(Respiratory)/Proning Guideline.pdf',
    'GI_Liver_and_Transplant/Prokinetics in ICU.pdf',
    'Drugs/Octreotide.pdf',
    'Procedures/ACCP CVC placement following completion of initial competencies.pdf',
    'Procedures/ACCPs acquiring initial CVC competencies.pdf',
    'Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf',
    'Delirium/Violence and Agression.pdf',
    'Post_op_care/Care of the Transgender Patient.pdf',
    'GI_Liver_and_Transplant/Plasma exchange in Acute Liver Failure.pdf',
    'End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf',
    'Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf',
    'Transfer/Transfer Guidelines.pdf',
    'Breathing(Respiratory)/Equipment/NIV through Nihon Kohden  Setup.pdf',
    'Infection_and_sepsis/Influenza in Critical Care.pdf',
    'Policies_and_admin/Anticipated Post op flow surgical patients.pdf',
    'Policies_and_admin/Discharge Home from Critical Care.pdf',
    'Policies_and_admin/Repatriaiton Checklist for Critical Care.pdf'
)

for ($i = 0; $i -lt $files.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $files[$i]
}
